# Apply roster update: reorder/replace player rows (A: Name, B: Positions, C: Team)
# Data rows run from row 2 through row 18 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Jared McCain",      "PG,SG",    "Philadelphia 76ers"),
    @("LaMelo Ball",       "PG,SG",    "Charlotte Hornets"),
    @("Brandon Miller",    "SG,SF",    "Charlotte Hornets"),
    @("Brandon Ingram",    "SG,SF,PF", "New Orleans Pelicans"),
    @("Bilal Coulibaly",   "SG,SF",    "Washington Wizards"),
    @("Julius Randle",     "PF",       "Minnesota Timberwolves"),
    @("Harrison Barnes",   "SF,PF",    "San Antonio Spurs"),
    @("Cameron Johnson",   "SF,PF",    "Brooklyn Nets"),
    @("Bam Adebayo",       "C",        "Miami Heat"),
    @("Anthony Davis",     "PF,C",     "Los Angeles Lakers"),
    @("Isaiah Hartenstein","C",        "Oklahoma City Thunder"),
    @("Ty Jerome",         "PG,SG",    "Cleveland Cavaliers"),
    @("Damian Lillard",    "PG",       "Milwaukee Bucks"),
    @("Derrick White",     "PG,SG",    "Boston Celtics"),
    @("Zion Williamson",   "PF,C",     "New Orleans Pelicans"),
    @("Cade Cunningham",   "PG,SG",    "Detroit Pistons"),
    @("Cam Thomas",        "SG,SF",    "Brooklyn Nets")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
